$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 12 (Leve Item ID 5515)
$ws.Range("H12").Value = 424.66666
$ws.Range("I12").Value = 424.66666
$ws.Range("K12").Value = 424.66666
$ws.Range("M12").Value = -254.66666

# Row 51 (Leve Item ID 5486)
$ws.Range("H51").Value = 199999
$ws.Range("J51").Value = 199999
$ws.Range("L51").Value = 199999
$ws.Range("N51").Value = -200967

# Row 64 (Leve Item ID 5506)
$ws.Range("H64").Value = 4985.294
$ws.Range("I64").Value = 5350
$ws.Range("K64").Value = 5350
$ws.Range("M64").Value = -5102

# Row 67 (Leve Item ID 5506)
$ws.Range("H67").Value = 4985.294
$ws.Range("I67").Value = 5350
$ws.Range("K67").Value = 5350
$ws.Range("M67").Value = -4492

# Row 75 (Leve Item ID 10654)
$ws.Range("H75").Value = 44000
$ws.Range("J75").Value = 60000
$ws.Range("L75").Value = 60000
$ws.Range("N75").Value = -61872

# Row 78 (Leve Item ID 10654)
$ws.Range("H78").Value = 44000
$ws.Range("J78").Value = 60000
$ws.Range("L78").Value = 180000
$ws.Range("N78").Value = -189360

# Row 87 (Leve Item ID 10651)
$ws.Range("H87").Value = 124970
$ws.Range("J87").Value = 124970
$ws.Range("L87").Value = 124970
$ws.Range("N87").Value = -127466

# Row 90 (Leve Item ID 10651)
$ws.Range("H90").Value = 124970
$ws.Range("J90").Value = 124970
$ws.Range("L90").Value = 374910
$ws.Range("N90").Value = -387390

# Row 106 (Leve Item ID 19903)
$ws.Range("H106").Value = 3290
$ws.Range("I106").Value = 1270.5714
$ws.Range("K106").Value = 1270.5714
$ws.Range("M106").Value = -639.5714

# Row 129 (Leve Item ID 36115)
$ws.Range("H129").Value = 764.5
$ws.Range("I129").Value = 557
$ws.Range("K129").Value = 1671
$ws.Range("M129").Value = 3329

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 4254.5835
$ws.Range("I137").Value = 1664.4375
$ws.Range("K137").Value = 4993.3125
$ws.Range("M137").Value = -2443.3125

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 34 (Leve Item ID 2753)
$ws.Range("H34").Value = 500000
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

# Row 46 (Leve Item ID 3498)
$ws.Range("H46").Value = 4898
$ws.Range("J46").Value = 4898
$ws.Range("L46").Value = 4898
$ws.Range("N46").Value = -5536

# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 10427279
$ws.Range("I74").Value = 17859516
$ws.Range("J74").Value = 22146.8
$ws.Range("K74").Value = 17859516
$ws.Range("L74").Value = 22146.8
$ws.Range("M74").Value = -17858642
$ws.Range("N74").Value = -23894.8

# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 10427279
$ws.Range("I77").Value = 17859516
$ws.Range("J77").Value = 22146.8
$ws.Range("K77").Value = 89297580
$ws.Range("L77").Value = 110734
$ws.Range("M77").Value = -89293212
$ws.Range("N77").Value = -119470

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 7017.5
$ws.Range("I132").Value = 4864.9556
$ws.Range("K132").Value = 14594.8668
$ws.Range("M132").Value = -12064.8668

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 88 (Leve Item ID 10626)
$ws.Range("H88").Value = 57719.5
$ws.Range("J88").Value = 57719.5
$ws.Range("L88").Value = 57719.5
$ws.Range("N88").Value = -58531.5

# Row 91 (Leve Item ID 10626)
$ws.Range("H91").Value = 57719.5
$ws.Range("J91").Value = 57719.5
$ws.Range("L91").Value = 57719.5
$ws.Range("N91").Value = -60527.5

# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 2380.7693
$ws.Range("I99").Value = 1922.7273
$ws.Range("K99").Value = 1922.7273
$ws.Range("M99").Value = -424.7273

# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 2396
$ws.Range("I105").Value = 1920.7778
$ws.Range("J105").Value = 2930.625
$ws.Range("K105").Value = 1920.7778
$ws.Range("L105").Value = 2930.625
$ws.Range("M105").Value = -173.7778000000001
$ws.Range("N105").Value = -6424.625

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 16 (Leve Item ID 27691)
$ws.Range("H16").Value = 2462.25
$ws.Range("I16").Value = 2462.25
$ws.Range("K16").Value = 2462.25
$ws.Range("M16").Value = -2175.25

# Row 100 (Leve Item ID 34388)
$ws.Range("H100").Value = 60993.332
$ws.Range("J100").Value = 60993.332
$ws.Range("L100").Value = 60993.332
$ws.Range("N100").Value = -63157.332

# Row 113 (Leve Item ID 27691)
$ws.Range("H113").Value = 2462.25
$ws.Range("I113").Value = 2462.25
$ws.Range("K113").Value = 2462.25
$ws.Range("M113").Value = -292.25

# Row 122 (Leve Item ID 36196)
$ws.Range("H122").Value = 1165
$ws.Range("I122").Value = 1153.2
$ws.Range("K122").Value = 3459.6
$ws.Range("M122").Value = -1009.6

# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 2337.5264
$ws.Range("I132").Value = 2356.3333
$ws.Range("J132").Value = 1999
$ws.Range("K132").Value = 7068.999899999999
$ws.Range("L132").Value = 5997
$ws.Range("M132").Value = -4538.999899999999
$ws.Range("N132").Value = -11057

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 11 (Leve Item ID 4745)
$ws.Range("H11").Value = 1873.5555
$ws.Range("I11").Value = 1969.48
$ws.Range("K11").Value = 5908.440000000001
$ws.Range("M11").Value = -5768.440000000001

# Row 39 (Leve Item ID 4712)
$ws.Range("H39").Value = 334335.34
$ws.Range("I39").Value = 1003
$ws.Range("J39").Value = 501001.5
$ws.Range("K39").Value = 3009
$ws.Range("L39").Value = 1503004.5
$ws.Range("M39").Value = -2715
$ws.Range("N39").Value = -1503592.5

# Row 46 (Leve Item ID 4701)
$ws.Range("H46").Value = 1107.3334
$ws.Range("J46").Value = 1656.1666
$ws.Range("L46").Value = 4968.4998
$ws.Range("N46").Value = -5150.4998

# Row 51 (Leve Item ID 4646)
$ws.Range("H51").Value = 15791.546
$ws.Range("I51").Value = 8398.857
$ws.Range("K51").Value = 25196.571
$ws.Range("M51").Value = -24736.571

# Row 102 (Leve Item ID 19813)
$ws.Range("H102").Value = 5000
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 0
$ws.Range("M102").Value = 15000
$ws.Range("N102").Value = -19868

# Row 122 (Leve Item ID 36078)
$ws.Range("H122").Value = 1838.75
$ws.Range("J122").Value = 2158.125
$ws.Range("L122").Value = 19423.125
$ws.Range("N122").Value = -24323.125

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 7392.5713
$ws.Range("I102").Value = 4458
$ws.Range("K102").Value = 4458
$ws.Range("M102").Value = -2836

# Row 109 (Leve Item ID 25691)
$ws.Range("H109").Value = 47748.5
$ws.Range("J109").Value = 47748.5
$ws.Range("L109").Value = 47748.5
$ws.Range("N109").Value = -49828.5

# Row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 3639.9524
$ws.Range("J113").Value = 3563.4285
$ws.Range("L113").Value = 3563.4285
$ws.Range("N113").Value = -7903.4285

# Row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 4407
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 4407
$ws.Range("K126").Value = 0
$ws.Range("M126").Value = 13221
$ws.Range("N126").Value = -18161

# Row 128 (Leve Item ID 34544)
$ws.Range("H128").Value = 84664.44500000001
$ws.Range("J128").Value = 84664.44500000001
$ws.Range("L128").Value = 84664.44500000001
$ws.Range("N128").Value = -94624.44500000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 127 (Leve Item ID 34401)
$ws.Range("H127").Value = 88000.8
$ws.Range("J127").Value = 88000.8
$ws.Range("L127").Value = 88000.8
$ws.Range("N127").Value = -97920.8

# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 47455.074
$ws.Range("I136").Value = 5015.4
$ws.Range("K136").Value = 15046.2
$ws.Range("M136").Value = -12496.2

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 59 (Leve Item ID 3201)
$ws.Range("H59").Value = 18426.666
$ws.Range("J59").Value = 18426.666
$ws.Range("L59").Value = 18426.666
$ws.Range("N59").Value = -19902.666

# Row 68 (Leve Item ID 10762)
$ws.Range("H68").Value = 60000
$ws.Range("J68").Value = 60000
$ws.Range("L68").Value = 60000
$ws.Range("N68").Value = -61622

# Row 71 (Leve Item ID 10762)
$ws.Range("H71").Value = 60000
$ws.Range("J71").Value = 60000
$ws.Range("L71").Value = 180000
$ws.Range("N71").Value = -188112

# Row 81 (Leve Item ID 12596)
$ws.Range("H81").Value = 3174.5
$ws.Range("I81").Value = 898.6667
$ws.Range("J81").Value = 10002
$ws.Range("K81").Value = 1797.3334
$ws.Range("L81").Value = 20004
$ws.Range("M81").Value = -736.3334
$ws.Range("N81").Value = -22126

# Row 84 (Leve Item ID 12596)
$ws.Range("H84").Value = 3174.5
$ws.Range("I84").Value = 898.6667
$ws.Range("J84").Value = 10002
$ws.Range("K84").Value = 8986.666999999999
$ws.Range("L84").Value = 100020
$ws.Range("M84").Value = -3682.666999999999
$ws.Range("N84").Value = -110628

# Row 100 (Leve Item ID 19981)
$ws.Range("H100").Value = 1013.5333
$ws.Range("I100").Value = 1014.5
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 2029
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -1488
$ws.Range("N100").Value = -3082

# Row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 20001940
$ws.Range("J107").Value = 1963.125
$ws.Range("L107").Value = 5889.375
$ws.Range("N107").Value = -9729.375

# Row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 2936
$ws.Range("I126").Value = 2010.8
$ws.Range("K126").Value = 6032.4
$ws.Range("M126").Value = -3562.4

# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 16195.192
$ws.Range("I136").Value = 777.3077
$ws.Range("J136").Value = 777.3077
$ws.Range("K136").Value = 2331.9231
$ws.Range("M136").Value = 218.0769

Write-Host "All updates applied"
